$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for range C2:O13 (columns C..O, rows 2..13)
$data = @(
    @(3.91, 20823.46, 0.1, 7.31, 91.22, 10468.56, 2475.96, 3.31, 15371.08, 1.37, 17022.78, 2.3, 66.46),
    @(2.29, 20568, 0.09, 6.22, 51.29, 10666.53, 2146.93, 2.92, 13890.11, 0.91, 17593.45, 2.06, 63.54),
    @(3.66, 9827.54, 0.07, 3.63, 107.58, 2008.12, 960.84, 1.64, 5537.32, 1.23, 5626.08, 0.84, 25.55),
    @(0.936, 0.472, 0.7, 0.497, 1.179, 0.192, 0.388, 0.495, 0.36, 0.898, 0.331, 0.365, 0.384),
    @(0.7, 4694.15, 0.03, 2.96, 14.12, 7046.92, 1156.24, 0.95, 7809.54, 0.34, 7209.95, 1.34, 25.79),
    @(12.66, 42456.22, 0.33, 16.97, 340.49, 13707.61, 4539.49, 7, 24318.26, 4.75, 28697.95, 4.03, 129.28),
    @(4.68, 19386.98, 0.09, 4.83, 40.03, 9585.25, 1858.3, 2.36, 9344.58, 1.17, 15448.08, 2.88, 48.85),
    @(2.86, 18114, 0.09, 4.14, 29.37, 8449.39, 1738.99, 1.8, 8562.15, 1.1, 15439.49, 2.7, 44.4),
    @(4.66, 7561.5, 0.03, 2.37, 36.28, 3312.65, 567.82, 2.33, 4330.96, 0.39, 4600.42, 0.83, 23.47),
    @(0.996, 0.39, 0.333, 0.491, 0.906, 0.346, 0.306, 0.987, 0.463, 0.333, 0.298, 0.288, 0.48),
    @(1.49, 9285.73, 0.05, 2.34, 8.61, 5430.34, 1094.05, 0.68, 4211.35, 0.56, 8841.53, 1.63, 19.63),
    @(16.15, 38408.84, 0.14, 12.13, 137.21, 16937.76, 3221.55, 8.89, 22339.51, 2.04, 26815.01, 4.58, 102.41)
)

$startRow = 2
$startCol = 3  # column C

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
